# Weekly refresh of the "Poroto granado" hortaliza data: the per-record
# fields (Fecha, Volumen, Precio minimo/maximo/promedio ponderado, Origen,
# Precio $/Kg) are reshuffled across the existing data rows (2-15, 17-19;
# row 16 is left untouched), while the record-invariant columns (Mercado,
# Region, Codreg, Categoria, Variedad, Calidad, Unidad de comercializacion,
# Kg o Unidades, Clasificacion) stay put on each row.
#
# destination row -> source row it now takes its Fecha/Volumen/Precio/Origen from
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 4
    3 = 8
    4 = 13
    5 = 17
    6 = 7
    7 = 9
    8 = 19
    9 = 14
    10 = 5
    11 = 15
    12 = 10
    13 = 6
    14 = 18
    15 = 3
    17 = 11
    18 = 12
    19 = 2
}

# Columns that move with the permutation:
# D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, O=Origen, P=Precio $/Kg
$cols = @(4, 10, 11, 12, 13, 15, 16)

$allRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 17, 18, 19)

# Snapshot every affected row's values BEFORE any writes, since this is a
# permutation (several destinations read from rows that are themselves
# about to be overwritten).
$snap = @{}
foreach ($r in $allRows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snap[$r] = $rowVals
}

# Write each destination row's values from its mapped source row's snapshot.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $snap[$srcRow][$c]
    }
}
